$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4, column B: "Heatshrink" -> "Heatsink"
$ws.Range("B4").Value = "Heatsink"

# Row 4, column H: add hyperlink to the CUI product page (text already shows the URL)
$ws.Hyperlinks.Add($ws.Range("H4"), "https://www.mouser.com/ProductDetail/CUI/HSE-B20250-040H?qs=sGAEpiMZZMttgyDkZ5Wiut%252B4GcHIZ2pKOgousR6bMSo%3D")
$ws.Range("H4").Style = "Hyperlink"

# Update the active selection to H4 (matches the saved selection state)
$null = $ws.Range("H4").Select()
